$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New header row (CIN/NOM/PRENOM/DEPARTEMENT -> CNE/LastName/...) ---
$ws.Range("A1").Value = "CNE"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "FirstName"
$ws.Range("D1").Value = "Dept_Attachement"
$ws.Range("E1").Value = "Phone"
$ws.Range("F1").Value = "Email"

# --- 2) Fix a name typo on the last existing row (Azzouzi -> El Azzouzi) ---
$ws.Range("B7").Value = "El Azzouzi"

# --- 3) Fill in the new Phone column ---------------------------------------
$ws.Range("E2").Value = "06 82 02 22 11"
$ws.Range("E3").Value = "02 51 40 28 06"
$ws.Range("E4").Value = "08 55 45 51 46"
$ws.Range("E5").Value = "03 42 67 10 87"
$ws.Range("E6").Value = "08 89 02 36 34"
$ws.Range("E7").Value = "08 44 72 02 14"

# --- 4) Fill in the new Email column -----------------------------------
$ws.Range("F2").Value = "mattis.Integer@tellus.net"
$ws.Range("F3").Value = "hendrerit.a.arcu@massanonante.org"
$ws.Range("F4").Value = "eros.non.enim@erosturpis.org"
$ws.Range("F5").Value = "mauris.sapien.cursus@dolor.net"
$ws.Range("F6").Value = "quis.turpis@congue.ca"
$ws.Range("F7").Value = "at.fringilla@etcommodoat.net"

# --- 5) Fill in the Dept_Attachement (department code) column --------------
$ws.Range("D2").Value = "SIC"
$ws.Range("D3").Value = "SIC"
$ws.Range("D4").Value = "MI"
$ws.Range("D5").Value = "MI"
$ws.Range("D6").Value = "SIC"
$ws.Range("D7").Value = "SIC"

# --- 6) Fill in the new CNE (row id) column ---------------------------------
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# --- 7) Phone/Email came in with their own (pasted-in) font: plain black,
#     no theme color - Excel records this as a second cell style. -----------
$normal2 = $wb.Styles.Add("Normal 2")
$normal2.Font.Color = 0

$ws.Range("E2:E7").Style = "Normal 2"
$ws.Range("F2:F8").Style = "Normal 2"

# --- 8) New Email column width ----------------------------------------
$ws.Columns("F").ColumnWidth = 21.16666666666667

# --- 9) Selection left where the user ended up editing ----------------------
$ws.Range("E9").Select() | Out-Null
